$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 (current row 4 becomes row 5),
# so that the existing row4 data shifts down and we can fill row4 with
# the new "entropy" row content after updating row2/row3.
$ws.Rows.Item(4).Insert()

# Row 2: A2 0 -> 29 ; C2 1 -> 0.96 (B2/D2 unchanged)
$ws.Range("A2").Value = 29
$ws.Range("C2").Value = 0.96

# Row 3: A3 3 -> 33 ; B3 text -> entropy ; C3 1 -> 0.88
$ws.Range("A3").Value = 33
$ws.Range("B3").Value = "{thumb, index} (distance) {diff} - |spectral| entropy"
$ws.Range("C3").Value = 0.88

# Row 4 (new row, inherits formatting from row above on insert): fill with
# the "rel_pwr_2_to_4" data that used to be row 3
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 35
$ws.Range("B4").Value = "{thumb, index} (distance) {diff} - |spectral| rel_pwr_2_to_4"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0

# Row 5 (previously row 4, shifted down): update A5 value
$ws.Range("A5").Value = 37
